# Applies the "Deploying to gh-pages" metadata refresh to the
# ValueSet-tnm-distant-metastases-category-vs workbook:
#   - rename the two "Include from SNOMED CT[ 2]" sheets to "Include #0"/"Include #1"
#   - bump the Metadata!Date value
#   - insert a new "Jurisdiction" property row right after "Contact" on the Metadata sheet

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Include from SNOMED CT" sheets -------------------------
$wsInclude0 = $wb.Worksheets.Item("Include from SNOMED CT")
$wsInclude0.Name = "Include #0"

$wsInclude1 = $wb.Worksheets.Item("Include from SNOMED CT 2")
$wsInclude1.Name = "Include #1"

# --- 2. Update the Metadata sheet ------------------------------------------
$ws = $wb.Worksheets.Item("Metadata")

# Refresh the "Date" property value (row 8: A8=Date, B8=<timestamp>)
$ws.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# Insert a new "Jurisdiction" property row directly under "Contact" (row 10),
# pushing Description/Purpose/Copyright/Immutable down by one row.
$ws.Rows.Item(11).Insert()

# Match the look of the surrounding property rows (border/alignment style).
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
